# Add a new address-book entry for "Angeline L. Rayray" at the bottom of the
# "Address" sheet (new row 33), following the same Name/Address pattern used
# by all the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row right after the last existing entry (row 32).
$ws.Range("A33").Value = "Angeline L. Rayray"
$ws.Range("B33").Value = "1735 Road 3, Bagong Sikat, Punta St., Sta. Ana Manila City"

# Reflect the cell the user ended up clicking on after typing the new entry.
$ws.Range("D31").Select()
